# Checkpoint for the evening
# Re-type a few existing cells with a leading apostrophe (forcing Excel's
# "Text" quote-prefix format and dropping the stray trailing apostrophe
# that had been baked into the old values), then append two new test rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: fix the quote-prefixed text columns (drop trailing ' artifact) ----
$ws.Range("B2").Value = "'QA Testing 091624/Cart 1 Indoor Data/"
$ws.Range("C2").Value = "'USB-TEMP (Device 0) - Analog - 9-16-2024 2-09-36.6 PM.csv"
$ws.Range("D2").Value = "'USB-1208FS-Plus (Device 1) - Analog - 9-16-2024 2-09-37.1 PM.csv"
$ws.Range("E2").Value = "'None"
$ws.Range("J2").Value = "'Cart 1 Temp"
$ws.Range("K2").Value = "'Cart 1 IR"

# ---- Row 3: same cleanup (B3 previously had no quote-prefix at all) ----
$ws.Range("B3").Value = "'QA Testing 091624/Cart 1 Outdoor Data/"
$ws.Range("C3").Value = "'USB-TEMP (Device 0) - Analog - 9-16-2024 1-09-41.2 PM.csv"
$ws.Range("D3").Value = "'USB-1208FS-Plus (Device 1) - Analog - 9-16-2024 1-09-41.7 PM.csv"
$ws.Range("E3").Value = "'None"
$ws.Range("J3").Value = "'Cart 1 Temp"
$ws.Range("K3").Value = "'Cart 1 IR"

# ---- Row 4: new test #3 (Cart 2 indoor QC pre-check) ----
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "'QA Testing 091624/Cart 2 Indoor Data/"
$ws.Range("C4").Value = "USB-TEMP (Device 0) - Analog - 9-16-2024 10-49-41.1 AM.csv"
$ws.Range("D4").Value = "USB-1208FS-Plus (Device 1) - Analog - 9-16-2024 10-49-40.838 AM.csv"
$ws.Range("E4").Value = "None"
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").Value = "9/16/2024"
$ws.Range("G4").Value = "indoor  QC pre-check"
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = "4-wire"
$ws.Range("J4").Value = "Cart 2 Temp"
$ws.Range("K4").Value = "Cart 2 IR"

# ---- Row 5: new test #4 (GPS / outdoor transect test) ----
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "GPS Test"
$ws.Range("C5").Value = "USB-TEMP (Device 1) - Analog - 8-9-2024 10-24-07.5 AM.csv"
$ws.Range("D5").Value = "USB-1208FS-Plus (Device 0) - Analog - 8-9-2024 10-24-08.0 AM.csv"
$ws.Range("E5").Value = "Midway080924east.csv"
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null
$ws.Range("F5").Value = "8/9/2024"
$ws.Range("G5").Value = "outdoor transect"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = "4-wire"
$ws.Range("J5").Value = "Cart 1 Temp"
$ws.Range("K5").Value = "Cart 1 IR"
